$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "['Weapon', 'Knife', 'Traditional weapons', 'knobkieries', 'kieries', 'spears', 'panga']"
$ws.Range("E2").Value = "['weapon,knife,traditional', 'weapons,knobkieries,kieries,spears,panga']"

$ws.Range("C3").Value = "['Throw', 'Threw', 'Attack', 'attacked ']"
$ws.Range("D3").Value = "['Stones', 'rocks', 'bricks', 'brike ']"
$ws.Range("E3").Value = "['throw,threw,attack,attacked']"
$ws.Range("F3").Value = "['stones,rocks,bricks,brike']"

$ws.Range("C4").Value = "['Gun', 'cannon', 'canon', 'shotgun', 's/gun', 'rubber', 'stungrenade', 'stun grenade', 'stuntgrenade', 'tear gas', 'teargas', 'bullets', 'rubberbullets']"
$ws.Range("E4").Value = "['gun,cannon,canon,shotgun,s/gun,rubber,stungrenade,stun', 'grenade,stuntgrenade,tear', 'gas,teargas,bullets,rubberbullets']"
